# Update "想去人数" (want-to-go count) values in sheet "展览" and "全部类型"
# to reflect newly generated output (gh-pages rebuild).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 281
$ws1.Range("F4").Value = 7741
$ws1.Range("F5").Value = 5636
$ws1.Range("F6").Value = 469
$ws1.Range("F10").Value = 255
$ws1.Range("F11").Value = 236
$ws1.Range("F12").Value = 56

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 281
$ws4.Range("F4").Value = 7741
$ws4.Range("F5").Value = 5636
$ws4.Range("F6").Value = 469
$ws4.Range("F10").Value = 255
$ws4.Range("F13").Value = 236
$ws4.Range("F14").Value = 56

$wb.Save()
